# Auto-generated edit script: applies the targeted cell-level corrections
# (idno swaps, lecturer-name corrections, and date/time corrections)
# described by the commit, matching the original cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain value corrections (numbers, times, names) ---
$ws.Range("C15").Value = 9521
$ws.Range("N15").Value = "Vaverka Lukáš, Mgr."
$ws.Range("C16").Value = 3528
$ws.Range("N16").Value = "Kučerová Silvie Rita, RNDr. Ph.D."
$ws.Range("H17").Value = "09:00"
$ws.Range("I17").Value = "12:50"
$ws.Range("H18").Value = "12:00"
$ws.Range("I18").Value = "15:50"
$ws.Range("C57").Value = 2855
$ws.Range("N57").Value = "Jirsák Jan, doc. RNDr. Ph.D."
$ws.Range("C58").Value = 2856
$ws.Range("N58").Value = "Škvorová Magda, Ing. Ph.D."
$ws.Range("C61").Value = 1064
$ws.Range("C62").Value = 1503
$ws.Range("C64").Value = 4310
$ws.Range("C65").Value = 4587
$ws.Range("N65").Value = "Ryšánek Petr, RNDr. Ph.D."
$ws.Range("C66").Value = 5431
$ws.Range("N66").Value = "Čermák Jan, doc. Ing. CSc."
$ws.Range("C67").Value = 2856
$ws.Range("N67").Value = "Škvorová Magda, Ing. Ph.D."
$ws.Range("N79").Value = "Henych Jiří, Ing. Ph.D."
$ws.Range("H81").Value = "15:00"
$ws.Range("I81").Value = "18:50"
$ws.Range("H82").Value = "09:00"
$ws.Range("I82").Value = "12:50"
$ws.Range("N87").Value = "Šícha Václav, RNDr. Ph.D."
$ws.Range("C105").Value = 4472
$ws.Range("N105").Value = "Liegertová Michaela, Mgr. Ph.D."
$ws.Range("C106").Value = 930
$ws.Range("N106").Value = "Ipser Jan, RNDr. CSc."
$ws.Range("C111").Value = 609
$ws.Range("C112").Value = 7640
$ws.Range("N114").Value = "Dočkal Jan, RNDr. Ph.D."
$ws.Range("N115").Value = "Lísal Martin, prof. Ing. DSc."

# --- Date-text corrections (datumOd/datumDo columns) ---
# These must stay as literal text like the original file (e.g. "10.12.2023"),
# not get auto-converted into a date serial number. Force text entry by
# temporarily marking the cell as Text, then restore the original look
# (General format, vertically centered) so the style matches the source.
$dateCells = @(
    @{Cell="G7"; Value="10.12.2023"},
    @{Cell="G71"; Value="10.12.2023"},
    @{Cell="F82"; Value="7.10.2023"},
    @{Cell="G82"; Value="7.10.2023"},
    @{Cell="F83"; Value="21.10.2023"},
    @{Cell="G83"; Value="21.10.2023"},
    @{Cell="F84"; Value="7.10.2023"},
    @{Cell="G84"; Value="7.10.2023"},
    @{Cell="F85"; Value="11.11.2023"},
    @{Cell="G85"; Value="11.11.2023"},
    @{Cell="G110"; Value="10.12.2023"}
)
foreach ($d in $dateCells) {
    $rng = $ws.Range($d.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $d.Value
    $rng.ClearFormats()
    $rng.VerticalAlignment = -4108
}
